$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("B2", 0.2073578595317726),
    @("C2", 0.5518394648829431),
    @("J2", 0.01003344481605351),
    @("P2", 0.1438127090301003),
    @("S2", 0.08695652173913043),
    @("B3", 0.005847953216374269),
    @("C3", 0.04678362573099415),
    @("J3", 0.02923976608187134),
    @("P3", 0.7543859649122807),
    @("S3", 0.1637426900584795),
    @("B6", 0.0532319391634981),
    @("D6", 0.007604562737642586),
    @("F6", 0.06083650190114068),
    @("J6", 0.2167300380228137),
    @("O6", 0.01901140684410646),
    @("Q6", 0.1596958174904943),
    @("R6", 0.05703422053231939),
    @("S6", 0.4258555133079848),
    @("B7", 0.09051724137931035),
    @("D7", 0.02586206896551724),
    @("E7", 0.004310344827586207),
    @("F7", 0.05603448275862069),
    @("J7", 0.08620689655172414),
    @("O7", 0.03017241379310345),
    @("Q7", 0.1939655172413793),
    @("R7", 0.06465517241379311),
    @("S7", 0.4482758620689655),
    @("B8", 0.1003649635036496),
    @("D8", 0.01642335766423358),
    @("F8", 0.06386861313868614),
    @("J8", 0.1058394160583942),
    @("O8", 0.01642335766423358),
    @("Q8", 0.1751824817518248),
    @("R8", 0.08029197080291971),
    @("S8", 0.4416058394160584),
    @("B9", 0.07755102040816327),
    @("D9", 0.00816326530612245),
    @("F9", 0.05714285714285714),
    @("J9", 0.08163265306122448),
    @("O9", 0.0163265306122449),
    @("Q9", 0.2122448979591837),
    @("R9", 0.07755102040816327),
    @("S9", 0.4693877551020408),
    @("B10", 0.09672262190247802),
    @("D10", 0.026378896882494),
    @("E10", 0.001598721023181455),
    @("F10", 0.08713029576338929),
    @("J10", 0.1167066346922462),
    @("O10", 0.01278976818545164),
    @("Q10", 0.1998401278976819),
    @("R10", 0.06314948041566747),
    @("S10", 0.3956834532374101),
    @("G11", 0.1233333333333333),
    @("J11", 0.05666666666666666),
    @("K11", 0.1633333333333333),
    @("L11", 0.63),
    @("S11", 0.02666666666666667),
    @("G12", 0.8071065989847716),
    @("J12", 0.1269035532994924),
    @("K12", 0.005076142131979695),
    @("L12", 0.01522842639593909),
    @("S12", 0.04568527918781726),
    @("G13", 0.75),
    @("J13", 0.2083333333333333),
    @("S13", 0.04166666666666666),
    @("G14", 1),
    @("F15", 0.01694915254237288),
    @("H15", 0.1991525423728814),
    @("I15", 0.09322033898305085),
    @("J15", 0.2923728813559322),
    @("K15", 0.0423728813559322),
    @("M15", 0.02542372881355932),
    @("O15", 0.05084745762711865),
    @("S15", 0.2796610169491525),
    @("F16", 0.02427184466019417),
    @("H16", 0.2233009708737864),
    @("I16", 0.07766990291262135),
    @("J16", 0.325242718446602),
    @("K16", 0.07281553398058252),
    @("M16", 0.05825242718446602),
    @("O16", 0.07281553398058252),
    @("S16", 0.145631067961165),
    @("F17", 0.02282157676348548),
    @("H17", 0.1701244813278008),
    @("I17", 0.1037344398340249),
    @("J17", 0.3817427385892116),
    @("K17", 0.08713692946058091),
    @("M17", 0.02282157676348548),
    @("O17", 0.07883817427385892),
    @("S17", 0.1327800829875519),
    @("F18", 0.01724137931034483),
    @("H18", 0.2586206896551724),
    @("I18", 0.05747126436781609),
    @("J18", 0.3160919540229885),
    @("K18", 0.08620689655172414),
    @("M18", 0.02873563218390805),
    @("N18", 0.005747126436781609),
    @("O18", 0.09195402298850575),
    @("S18", 0.1379310344827586),
    @("F19", 0.01128818061088977),
    @("H19", 0.2177954847277556),
    @("I19", 0.100265604249668),
    @("J19", 0.351925630810093),
    @("K19", 0.1069057104913679),
    @("M19", 0.02523240371845949),
    @("O19", 0.06175298804780877),
    @("S19", 0.1248339973439575)
)

foreach ($change in $changes) {
    $ws.Range($change[0]).Value = $change[1]
}

Write-Output "Applied $($changes.Count) cell updates"